$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H129").Value = 899.16437
$ws.Range("I129").Value = 1219.8
$ws.Range("K129").Value = 3659.4
$ws.Range("M129").Value = 1340.6
$ws.Range("H132").Value = 1681.091
$ws.Range("I132").Value = 1681.091
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5043.272999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -2513.272999999999
$ws.Range("H137").Value = 3231.611
$ws.Range("I137").Value = 1270
$ws.Range("J137").Value = 4800.9
$ws.Range("K137").Value = 3810
$ws.Range("L137").Value = 14402.7
$ws.Range("M137").Value = -1260
$ws.Range("N137").Value = -19502.7
$ws.Range("H138").Value = 2548.0833
$ws.Range("J138").Value = 2829.6924
$ws.Range("L138").Value = 8489.0772
$ws.Range("N138").Value = -18769.0772

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3322810.5
$ws.Range("I2").Value = 3322810.5
$ws.Range("K2").Value = 3322810.5
$ws.Range("M2").Value = -3322697.5
$ws.Range("I45").Value = 1077.3
$ws.Range("J45").Value = 1808.2
$ws.Range("K45").Value = 1077.3
$ws.Range("L45").Value = 1808.2
$ws.Range("M45").Value = -700.3
$ws.Range("N45").Value = -2562.2
$ws.Range("H61").Value = 5014.2
$ws.Range("I61").Value = 2815.1538
$ws.Range("J61").Value = 9098.143
$ws.Range("K61").Value = 2815.1538
$ws.Range("L61").Value = 9098.143
$ws.Range("M61").Value = -2603.1538
$ws.Range("N61").Value = -9522.143
$ws.Range("H116").Value = 3322810.5
$ws.Range("I116").Value = 3322810.5
$ws.Range("K116").Value = 3322810.5
$ws.Range("M116").Value = -3320516.5
$ws.Range("H136").Value = 5014.2
$ws.Range("I136").Value = 2815.1538
$ws.Range("J136").Value = 9098.143
$ws.Range("K136").Value = 8445.4614
$ws.Range("L136").Value = 27294.429
$ws.Range("M136").Value = -5895.4614
$ws.Range("N136").Value = -32394.429

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3322810.5
$ws.Range("I3").Value = 3322810.5
$ws.Range("K3").Value = 3322810.5
$ws.Range("M3").Value = -3322696.5
$ws.Range("H22").Value = 533
$ws.Range("I22").Value = 533
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 533
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -360
$ws.Range("H134").Value = 10541.777
$ws.Range("I134").Value = 11271.2
$ws.Range("J134").Value = 8457.714
$ws.Range("K134").Value = 33813.60000000001
$ws.Range("L134").Value = 25373.142
$ws.Range("M134").Value = -31278.60000000001
$ws.Range("N134").Value = -30443.142

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2509.7837
$ws.Range("I31").Value = 1776.7931
$ws.Range("J31").Value = 5166.875
$ws.Range("K31").Value = 1776.7931
$ws.Range("L31").Value = 5166.875
$ws.Range("M31").Value = -1481.7931
$ws.Range("N31").Value = -5756.875
$ws.Range("H34").Value = 2509.7837
$ws.Range("I34").Value = 1776.7931
$ws.Range("J34").Value = 5166.875
$ws.Range("K34").Value = 1776.7931
$ws.Range("L34").Value = 5166.875
$ws.Range("M34").Value = -1574.7931
$ws.Range("N34").Value = -5570.875
$ws.Range("H74").Value = 27499.8
$ws.Range("J74").Value = 27499.8
$ws.Range("L74").Value = 27499.8
$ws.Range("N74").Value = -29247.8
$ws.Range("H77").Value = 27499.8
$ws.Range("J77").Value = 27499.8
$ws.Range("L77").Value = 82499.39999999999
$ws.Range("N77").Value = -91235.39999999999
$ws.Range("H122").Value = 1301.4193
$ws.Range("I122").Value = 1232.6666
$ws.Range("J122").Value = 1445.8
$ws.Range("K122").Value = 3697.9998
$ws.Range("L122").Value = 4337.4
$ws.Range("M122").Value = -1247.9998
$ws.Range("N122").Value = -9237.4
$ws.Range("H132").Value = 2514.3125
$ws.Range("I132").Value = 1201.4166
$ws.Range("J132").Value = 6453
$ws.Range("K132").Value = 3604.2498
$ws.Range("L132").Value = 19359
$ws.Range("M132").Value = -1074.2498
$ws.Range("N132").Value = -24419
$ws.Range("H134").Value = 855.1923
$ws.Range("I134").Value = 842.5417
$ws.Range("K134").Value = 2527.6251
$ws.Range("M134").Value = 7.374899999999798
$ws.Range("H140").Value = 14374.875
$ws.Range("J140").Value = 14374.875
$ws.Range("L140").Value = 14374.875
$ws.Range("N140").Value = -24734.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3717.5334
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("H109").Value = 4498
$ws.Range("I109").Value = 899.6
$ws.Range("K109").Value = 2698.8
$ws.Range("M109").Value = -1658.8
$ws.Range("H113").Value = 5925.5264
$ws.Range("J113").Value = 745.2
$ws.Range("L113").Value = 2235.6
$ws.Range("N113").Value = -6575.6
$ws.Range("H131").Value = 783.85
$ws.Range("J131").Value = 795.71576
$ws.Range("L131").Value = 2387.14728
$ws.Range("N131").Value = -12467.14728
$ws.Range("H140").Value = 1682.5151
$ws.Range("I140").Value = 867.8
$ws.Range("J140").Value = 2361.4443
$ws.Range("K140").Value = 2603.4
$ws.Range("L140").Value = 7084.3329
$ws.Range("M140").Value = 2576.6
$ws.Range("N140").Value = -17444.3329

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1751661.2
$ws.Range("I132").Value = 2749086.2
$ws.Range("J132").Value = 6167.75
$ws.Range("K132").Value = 8247258.600000001
$ws.Range("L132").Value = 18503.25
$ws.Range("M132").Value = -8244728.600000001
$ws.Range("N132").Value = -23563.25
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 31499.5
$ws.Range("J141").Value = 31499.5
$ws.Range("L141").Value = 31499.5
$ws.Range("N141").Value = -41859.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1902.2759
$ws.Range("I132").Value = 1765.2142
$ws.Range("J132").Value = 2030.2
$ws.Range("K132").Value = 5295.642599999999
$ws.Range("L132").Value = 6090.6
$ws.Range("M132").Value = -2765.642599999999
$ws.Range("N132").Value = -11150.6
$ws.Range("H136").Value = 3914.1428
$ws.Range("I136").Value = 2142.5715
$ws.Range("K136").Value = 6427.7145
$ws.Range("M136").Value = -3877.7145

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 1677075.4
$ws.Range("I63").Value = 6226
$ws.Range("J63").Value = 2512500
$ws.Range("K63").Value = 6226
$ws.Range("L63").Value = 2512500
$ws.Range("M63").Value = -5602
$ws.Range("N63").Value = -2513748
$ws.Range("H66").Value = 1677075.4
$ws.Range("I66").Value = 6226
$ws.Range("J66").Value = 2512500
$ws.Range("K66").Value = 18678
$ws.Range("L66").Value = 7537500
$ws.Range("M66").Value = -15558
$ws.Range("N66").Value = -7543740
$ws.Range("H103").Value = 25799.5
$ws.Range("J103").Value = 25799.5
$ws.Range("L103").Value = 25799.5
$ws.Range("N103").Value = -28143.5
$ws.Range("H132").Value = 2438.1052
$ws.Range("I132").Value = 2148.1538
$ws.Range("K132").Value = 6444.4614
$ws.Range("M132").Value = -3914.4614
$ws.Range("H136").Value = 18520070
$ws.Range("I136").Value = 26456050
$ws.Range("K136").Value = 79368150
$ws.Range("M136").Value = -79365600
$ws.Range("H140").Value = 56500
$ws.Range("J140").Value = 56500
$ws.Range("L140").Value = 56500
$ws.Range("N140").Value = -66860
$ws.Range("H141").Value = 73538.38
$ws.Range("J141").Value = 73538.38
$ws.Range("L141").Value = 73538.38
$ws.Range("N141").Value = -83898.38
